$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "08/08/2024 15:23"
$ws.Range("B8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0

$ws.Range("A9").Value = "08/08/2024 19:24"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 3753
$ws.Range("D9").Value = 4170.0
$ws.Range("E9").Value = 417.0
$ws.Range("F9").Value = 3753
$ws.Range("G9").Value = 0
